# Germanize the "Assignment13b" workbook: translate sheet/table text to German
# and switch the data columns' number format from USD to EUR accounting style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Sheet (tab) name -------------------------------------------------
$ws.Name = "2013 Verkaufszahlen (Detail)"

# --- Title cell ---------------------------------------------------------
$ws.Range("A1").Value = "Westbrook Parker Verkaufszahlen"

# --- Table header row (also renames the Excel Table's column headers) ---
$ws.Range("A2").Value = "Verkäufer"
$ws.Range("B2").Value = "Mai"
$ws.Range("C2").Value = "Juni"
$ws.Range("D2").Value = "Juli"
$ws.Range("E2").Value = "August"
$ws.Range("F2").Value = "September"
$ws.Range("G2").Value = "Oktober"

# --- Number format: US dollars -> Euro accounting format ----------------
$euroFormat = "_-* #,##0.00\ [`$€-407]_-;\-* #,##0.00\ [`$€-407]_-;_-* ""-""??\ [`$€-407]_-;_-@_-"
$ws.Range("B3:G32").NumberFormat = $euroFormat

# --- Restore the selected cell (was G18, now A9) -------------------------
$ws.Range("A9").Select()
